# Update "想去人数" (want-to-go count) figures on the 展览 (Exhibition) and
# 全部类型 (All types) sheets, refreshed as of the gh-pages data regeneration
# at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 549
$ws1.Range("F7").Value  = 1700
$ws1.Range("F10").Value = 31
$ws1.Range("F11").Value = 1711
$ws1.Range("F13").Value = 91
$ws1.Range("F19").Value = 35
$ws1.Range("F21").Value = 492
$ws1.Range("F23").Value = 163
$ws1.Range("F24").Value = 233
$ws1.Range("F25").Value = 252

# Sheet "全部类型" (All types) - aggregated view, one extra row vs 展览
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 549
$ws4.Range("F7").Value  = 1700
$ws4.Range("F11").Value = 31
$ws4.Range("F12").Value = 1711
$ws4.Range("F14").Value = 91
$ws4.Range("F20").Value = 35
$ws4.Range("F22").Value = 492
$ws4.Range("F24").Value = 163
$ws4.Range("F25").Value = 233
$ws4.Range("F26").Value = 252
